# Auto - Update data with bot!
# Update titles/links in the Blogs_used_list worksheet to reflect the
# latest scraped blog posts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: title update (same source/link stays, only the title text changes)
$ws.Range("D4").Value = "[PyTorch] RNN Layer 입출력 파라미터와 차원(shape) 이해"

# Row 5: new post - title + link
$ws.Range("D5").Value = "MATLAB으로 email 보내기"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2022/07/01/MATLAB_email.html"

# Row 26: title update
$ws.Range("D26").Value = "ai plus(est soft)"

# Row 36: new seminar post - title + link
$ws.Range("D36").Value = "Out-of-distribution Detection in image classification"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/369"

# Row 37: new seminar post - title + link
$ws.Range("D37").Value = "[Paper Review] Will I Sound Like Me? Improving Persona Consistency in Dialogues through Pragmatic Self-Consciousness"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=2102&mod=document&pageid=1"

# Row 49: title update
$ws.Range("D49").Value = "[NLP] 밑바닥부터 시작하는 딥러닝2 - Ch4 : word2vec 개선"

# Row 50: title update
$ws.Range("D50").Value = "공지예외주장 제도"

# Row 51: new post - title + link
$ws.Range("D51").Value = "[pandas] 데이터프레임 컬럼 순서 변경"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-%EB%8D%B0%EC%9D%B4%ED%84%B0%ED%94%84%EB%A0%88%EC%9E%84-%EC%BB%AC%EB%9F%BC-%EC%88%9C%EC%84%9C-%EB%B3%80%EA%B2%BD"
